$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.102.91"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "'1.652.69"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "'218.84"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'0.5254"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "'0.2662"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "'0.06349"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'20.59"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").Value = "'0.07706"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'4.607"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'1.670.49"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "'0.5604"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'0.0₅8190"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D18").Value = "'26.108.68"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'4.698"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'190.97"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").Value = "'5.985"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'145.88"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'0.1201"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").Value = "'7.255"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'15.92"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "'1.501"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'0.05629"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").Value = "'1.273"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "'3.379"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'1.580"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").Value = "'2.797"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'0.9477"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").Value = "'0.5758"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'5.973"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'0.8395"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'101.63"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "'1.014.44"
$ws.Range("E44").Value = "  -5.60%  "
$ws.Range("D45").Value = "'1.790.94"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "'58.23"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'0.05325"
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4344"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.000"
$ws.Range("E51").Value = "  -1.05%  "
